$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 15 (shifts rows 15-20 down to 16-21)
$ws.Rows.Item(15).Insert()

# Fill in the new row 15 with the new task
$ws.Cells.Item(15, 2).Value = 3.13
$ws.Cells.Item(15, 3).Value = "Xây dựng module thêm,sửa,xóa,tìm kiếm nhà xuất bản"
$ws.Cells.Item(15, 4).Value = "Đặng Quốc Đạt"
$ws.Cells.Item(15, 5).Value = "           27/10"

# Re-sequence the task numbers in column B for the shifted rows (16-19)
$ws.Cells.Item(16, 2).Value = 3.14
$ws.Cells.Item(17, 2).Value = 3.15
$ws.Cells.Item(18, 2).Value = 3.16
$ws.Cells.Item(19, 2).Value = 3.17

# Update the selection to match the target state
$ws.Range("D19").Select()
